$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.423.83'
$ws.Range('E2').Value = '  -0.34%  '
$ws.Range('D3').Value = '1.583.91'
$ws.Range('E3').Value = '  -0.20%  '
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('E5').Value = '  +0.42%  '
$ws.Range('E6').Value = '  -0.30%  '
$ws.Range('E7').Value = '  +0.11%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '44.46'
$ws.Range('E8').Value = '  +0.19%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '23.92'
$ws.Range('E9').Value = '  -1.51%  '
$ws.Range('E10').Value = '  -1.97%  '
$ws.Range('E11').Value = '  -1.58%  '
$ws.Range('E12').Value = '  +0.95%  '
$ws.Range('D13').Value = '1.810.57'
$ws.Range('E13').Value = '  -0.18%  '
$ws.Range('D14').Value = '1.586.57'
$ws.Range('E14').Value = '  -0.19%  '
$ws.Range('E15').Value = '  -0.77%  '
$ws.Range('E16').Value = '  -1.63%  '
$ws.Range('D17').Value = '28.459.88'
$ws.Range('E17').Value = '  -0.26%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '62.16'
$ws.Range('E18').Value = '  -1.40%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '230.08'
$ws.Range('E19').Value = '  -0.64%  '
$ws.Range('E20').Value = '  -0.73%  '
$ws.Range('D21').Value = '0.0₃0689'
$ws.Range('E21').Value = '  -2.42%  '
$ws.Range('E22').Value = '  +0.07%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '3.91'
$ws.Range('E23').Value = '  -3.13%  '
$ws.Range('E24').Value = '  -1.85%  '
$ws.Range('E25').Value = '  +3.33%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '151.95'
$ws.Range('E26').Value = '  +0.01%  '
$ws.Range('E27').Value = '  -1.42%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '6.43'
$ws.Range('E28').Value = '  -1.63%  '
$ws.Range('E31').Value = '  +2.85%  '
$ws.Range('E32').Value = '  -1.02%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '3.21'
$ws.Range('E33').Value = '  -1.26%  '
$ws.Range('E34').Value = '  -2.47%  '
$ws.Range('D35').Value = '1.396.14'
$ws.Range('E35').Value = '  +0.45%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '1.08'
$ws.Range('E36').Value = '  +7.00%  '
$ws.Range('E37').Value = '  -4.30%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '2.37'
$ws.Range('E38').Value = '  +0.55%  '
$ws.Range('E39').Value = '  +1.42%  '
$ws.Range('E40').Value = '  -0.75%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.522'
$ws.Range('E41').Value = '  -3.31%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.790'
$ws.Range('E43').Value = '  -2.60%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '1.88'
$ws.Range('E44').Value = '  +0.84%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.0461'
$ws.Range('E45').Value = '  -1.06%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '5.44'
$ws.Range('E46').Value = '  -3.32%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.961'
$ws.Range('E47').Value = '  -1.95%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '63.01'
$ws.Range('E48').Value = '  +0.10%  '
$ws.Range('D49').Value = '1.721.84'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '86.55'
$ws.Range('E50').Value = '  -0.61%  '
$ws.Range('D51').Value = '0.0₆0102'
$ws.Range('E51').Value = '  -2.60%  '
